$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Bump the cached "datetimeFigureOut" footer date from 11/12/2012
#    to 11/13/2012 on the slide master and every slide layout.
# ------------------------------------------------------------------
function Set-DatePlaceholderText($shapes, $newText) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes "11/13/2012"

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Set-DatePlaceholderText $layouts.Item($i).Shapes "11/13/2012"
}

# ------------------------------------------------------------------
# 2) Slide 2 ("Topics"): collapse the "ProtoDebugger" bullet's
#    " " + "GUI " runs into a single " GUI " run (text unchanged).
# ------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$contentShape = $slide2.Shapes.Item(2)
$tr2 = $contentShape.TextFrame.TextRange
$full2 = $tr2.Text
$marker = "ProtoDebugger GUI ScreenShots"
$pos = $full2.IndexOf($marker)
if ($pos -ge 0) {
    $start = $pos + 1 + ("ProtoDebugger").Length
    $len = " GUI ".Length
    $run = $tr2.Characters($start, $len)
    $run.Text = " GUI "
}

# ------------------------------------------------------------------
# 3) Slide 9 ("Lessons Learned..."): collapse the title's
#    "Lessons " + "Learned And Future Enhancements" runs into one.
# ------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$titleShape = $slide9.Shapes.Item(1)
$tr9 = $titleShape.TextFrame.TextRange
$full9 = $tr9.Text
$whole = $tr9.Characters(1, $full9.Length)
$whole.Text = $full9
